# "Add files via upload" commit — removes the hidden "Schema" slide (the
# 5th slide, show="0") and refreshes the cached "Update automatically"
# date placeholder text (datetimeFigureOut field) from 6/8/2020 to
# 6/12/2020 everywhere it is cached (every slide layout + the notes
# master), as happens when the deck is re-opened/re-saved on a later day.

$p = $ppt.ActivePresentation

# --- 1. Delete the 5th slide (hidden "Schema" slide, id 1205) ----------
# Removing it also drops its <p:sldId>/<p14:sldId> entries from the
# slide list and the "Introduction" section automatically.
$p.Slides.Item(5).Delete()

# --- 2. Refresh the cached date field on every slide layout -------------
$master = $p.Designs.Item(1).SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = "6/12/2020"
        }
    }
}

# --- 3. Refresh the cached date field on the notes master ---------------
$notesMaster = $p.NotesMaster
for ($si = 1; $si -le $notesMaster.Shapes.Count; $si++) {
    $shape = $notesMaster.Shapes.Item($si)
    if ($shape.Name -like "Date Placeholder*") {
        $shape.TextFrame.TextRange.Text = "6/12/2020"
    }
}

Write-Output "Slides remaining: $($p.Slides.Count)"
